$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 entirely (data now only spans rows 1-2)
$ws.Rows.Item(3).Delete() | Out-Null

# Row 1 values
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "Ab"
$ws.Range("C1").Value = "demoMed"
$ws.Range("D1").Value = 1.5
$ws.Range("E1").Value = 10
$ws.Range("F1").Value = 100
$ws.Range("G1").Value = "red"
$ws.Range("H1").Value = "shape"

# Row 2 values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Rx"
$ws.Range("C2").Value = "Paracetamol"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 225
$ws.Range("G2").Value = "đỏ"
$ws.Range("H2").Value = "hơi tròn"

# Update selection to match the new target cell
$ws.Range("G2").Select() | Out-Null
